$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9353491011476933
$ws.Range("D2").Value = 0.02470569264311706
$ws.Range("E2").Value = 0.3764878169001253
$ws.Range("F2").Value = 0.5320487359019879
$ws.Range("G2").Value = 0.3731684582329748
$ws.Range("H2").Value = 0.5359209923621577
$ws.Range("K2").Value = 0.3990874383953837
$ws.Range("L2").Value = 0.120836204172349
$ws.Range("N2").Value = 1.634832354216272
$ws.Range("O2").Value = 1.758896973849332
$ws.Range("B3").Value = 0.9094003381312916
$ws.Range("D3").Value = 0.02219258204841168
$ws.Range("E3").Value = 0.3798043307007166
$ws.Range("F3").Value = 0.5293819867776861
$ws.Range("G3").Value = 0.3719481265385625
$ws.Range("H3").Value = 0.5384930335440714
$ws.Range("K3").Value = 0.3556188903764053
$ws.Range("L3").Value = 0.1105801642356283
$ws.Range("N3").Value = 1.650987400577201
$ws.Range("O3").Value = 1.76146428061783
$ws.Range("B4").Value = 0.8938917477824759
$ws.Range("D4").Value = 0.02063884706763019
$ws.Range("E4").Value = 0.3820020767812533
$ws.Range("F4").Value = 0.5280712361764301
$ws.Range("G4").Value = 0.3714576501076721
$ws.Range("H4").Value = 0.5403000602223784
$ws.Range("K4").Value = 0.32876758247005
$ws.Range("L4").Value = 0.1043193871062016
$ws.Range("N4").Value = 1.661451150759669
$ws.Range("O4").Value = 1.764057927020474
$ws.Range("B5").Value = 0.887679092868467
$ws.Range("D5").Value = 0.02000303558871508
$ws.Range("E5").Value = 0.3829382973892432
$ws.Range("F5").Value = 0.5276192975582177
$ws.Range("G5").Value = 0.3713228427028952
$ws.Range("H5").Value = 0.5410937801045179
$ws.Range("K5").Value = 0.3177855593612406
$ws.Range("L5").Value = 0.1017773643608848
$ws.Range("N5").Value = 1.665852188541139
$ws.Range("O5").Value = 1.765370766779483
$ws.Range("B6").Value = 0.886653979508651
$ws.Range("D6").Value = 0.01989730036356718
$ws.Range("E6").Value = 0.3830962106363387
$ws.Range("F6").Value = 0.527549220270437
$ws.Range("G6").Value = 0.3713043874124509
$ws.Range("H6").Value = 0.5412290418746792
$ws.Range("K6").Value = 0.3159596136803913
$ws.Range("L6").Value = 0.1013558284303571
$ws.Range("N6").Value = 1.666591252217977
$ws.Range("O6").Value = 1.765604221876018
$ws.Range("B7").Value = 0.893807526828283
$ws.Range("D7").Value = 0.0206302829863958
$ws.Range("E7").Value = 0.382014538455504
$ws.Range("F7").Value = 0.528064808251159
$ws.Range("G7").Value = 0.3714555686188703
$ws.Range("H7").Value = 0.5403105323468083
$ws.Range("K7").Value = 0.3286196356812923
$ws.Range("L7").Value = 0.1042850667005695
$ws.Range("N7").Value = 1.661509950009377
$ws.Range("O7").Value = 1.764074596147751
$ws.Range("B8").Value = 0.9263142987080357
$ws.Range("D8").Value = 0.02384140648194943
$ws.Range("E8").Value = 0.377597876385126
$ws.Range("F8").Value = 0.5310614829185667
$ws.Range("G8").Value = 0.3726939611065205
$ws.Range("H8").Value = 0.536760593713943
$ws.Range("K8").Value = 0.3841333976056944
$ws.Range("L8").Value = 0.1172924225106584
$ws.Range("N8").Value = 1.640289642737628
$ws.Range("O8").Value = 1.759571092704419
$ws.Range("B9").Value = 0.9934028348092454
$ws.Range("D9").Value = 0.03005255372126925
$ws.Range("E9").Value = 0.3702157228014258
$ws.Range("F9").Value = 0.5395279157590522
$ws.Range("G9").Value = 0.3771774047130521
$ws.Range("H9").Value = 0.5316039469785352
$ws.Range("K9").Value = 0.4916888422650061
$ws.Range("L9").Value = 0.1430850901701319
$ws.Range("N9").Value = 1.602993787077679
$ws.Range("O9").Value = 1.758809453234903
$ws.Range("B10").Value = 1.044707327267048
$ws.Range("D10").Value = 0.03456238738617401
$ws.Range("E10").Value = 0.3655694244574512
$ws.Range("F10").Value = 0.547326260296785
$ws.Range("G10").Value = 0.3817271525695247
$ws.Range("H10").Value = 0.5289124567204198
$ws.Range("K10").Value = 0.5698854342697359
$ws.Range("L10").Value = 0.1622052788071784
$ws.Range("N10").Value = 1.578219280560994
$ws.Range("O10").Value = 1.763169548732321
$ws.Range("B11").Value = 1.068479906441723
$ws.Range("D11").Value = 0.03660218636001389
$ws.Range("E11").Value = 0.3636240204632024
$ws.Range("F11").Value = 0.5512164966591513
$ws.Range("G11").Value = 0.3840703282563425
$ws.Range("H11").Value = 0.5279256145779101
$ws.Range("K11").Value = 0.605274619536317
$ws.Range("L11").Value = 0.1709398732498499
$ws.Range("N11").Value = 1.56751802790691
$ws.Range("O11").Value = 1.76622155595868
$ws.Range("B12").Value = 1.077543838690104
$ws.Range("D12").Value = 0.03737288712032694
$ws.Range("E12").Value = 0.3629115003452164
$ws.Range("F12").Value = 0.55273886948887
$ws.Range("G12").Value = 0.3849969866631966
$ws.Range("H12").Value = 0.5275860229429981
$ws.Range("K12").Value = 0.6186486662097934
$ws.Range("L12").Value = 0.174252616442601
$ws.Range("N12").Value = 1.56354748253607
$ws.Range("O12").Value = 1.76753089037885
$ws.Range("B13").Value = 1.075589021291904
$ws.Range("D13").Value = 0.03720698032449832
$ws.Range("E13").Value = 0.3630638801241624
$ws.Range("F13").Value = 0.5524088109382745
$ws.Range("G13").Value = 0.3847956637909817
$ws.Range("H13").Value = 0.527657644126009
$ws.Range("K13").Value = 0.6157695422862446
$ws.Range("L13").Value = 0.1735389314612092
$ws.Range("N13").Value = 1.564398973092576
$ws.Range("O13").Value = 1.767242071033991
$ws.Range("B14").Value = 1.069224367162832
$ws.Range("D14").Value = 0.03666562720027144
$ws.Range("E14").Value = 0.3635649169172446
$ws.Range("F14").Value = 0.5513407569961402
$ws.Range("G14").Value = 0.3841457762649867
$ws.Range("H14").Value = 0.5278969929737656
$ws.Range("K14").Value = 0.6063754573790163
$ws.Range("L14").Value = 0.171212312374891
$ws.Range("N14").Value = 1.567189729642013
$ws.Range("O14").Value = 1.766326197527775
$ws.Range("B15").Value = 1.065333858820168
$ws.Range("D15").Value = 0.03633380677396758
$ws.Range("E15").Value = 0.3638749622239175
$ws.Range("F15").Value = 0.5506929518452353
$ws.Range("G15").Value = 0.3837528264242707
$ws.Range("H15").Value = 0.5280480408532924
$ws.Range("K15").Value = 0.6006177607566201
$ws.Range("L15").Value = 0.1697878564842341
$ws.Range("N15").Value = 1.568909799548859
$ws.Range("O15").Value = 1.765785200868635
$ws.Range("B16").Value = 1.043162411419758
$ws.Range("D16").Value = 0.03442884213365716
$ws.Range("E16").Value = 0.3656999438357929
$ws.Range("F16").Value = 0.5470789157595704
$ws.Range("G16").Value = 0.3815795255385552
$ws.Range("H16").Value = 0.5289817244034367
$ws.Range("K16").Value = 0.5675689221291975
$ws.Range("L16").Value = 0.1616351796360789
$ws.Range("N16").Value = 1.578930077688334
$ws.Range("O16").Value = 1.762991590782462
$ws.Range("B17").Value = 1.02967163793025
$ws.Range("D17").Value = 0.03325717261780881
$ws.Range("E17").Value = 0.3668625764690816
$ws.Range("F17").Value = 0.5449495612508528
$ws.Range("G17").Value = 0.3803163358790584
$ws.Range("H17").Value = 0.529615308859789
$ws.Range("K17").Value = 0.547247147566992
$ws.Range("L17").Value = 0.156643083753977
$ws.Range("N17").Value = 1.58522287121211
$ws.Range("O17").Value = 1.761551451768014
$ws.Range("B18").Value = 1.021952980644954
$ws.Range("D18").Value = 0.03258215588166991
$ws.Range("E18").Value = 0.3675471273669899
$ws.Range("F18").Value = 0.5437570758933958
$ws.Range("G18").Value = 0.3796155217886792
$ws.Range("H18").Value = 0.5300020904905267
$ws.Range("K18").Value = 0.5355414373732401
$ws.Range("L18").Value = 0.1537752273030577
$ws.Range("N18").Value = 1.588895857400363
$ws.Range("O18").Value = 1.760823704924292
$ws.Range("B19").Value = 1.019346615854147
$ws.Range("D19").Value = 0.0323534184435843
$ws.Range("E19").Value = 0.3677816251229835
$ws.Range("F19").Value = 0.5433588638391598
$ws.Range("G19").Value = 0.379382658173256
$ws.Range("H19").Value = 0.5301368901615575
$ws.Range("K19").Value = 0.5315751609907124
$ws.Range("L19").Value = 0.152804819742272
$ws.Range("N19").Value = 1.590148662781843
$ws.Range("O19").Value = 1.760594579208544
$ws.Range("B20").Value = 1.031103527504854
$ws.Range("D20").Value = 0.03338201331012414
$ws.Range("E20").Value = 0.3667371735419529
$ws.Range("F20").Value = 0.5451728959967852
$ws.Range("G20").Value = 0.380448140490472
$ws.Range("H20").Value = 0.5295455488534344
$ws.Range("K20").Value = 0.5494122170988192
$ws.Range("L20").Value = 0.1571741432676816
$ws.Range("N20").Value = 1.584547451683555
$ws.Range("O20").Value = 1.76169434744665
$ws.Range("B21").Value = 1.071092148846816
$ws.Range("D21").Value = 0.03682468290219276
$ws.Range("E21").Value = 0.3634170947326503
$ws.Range("F21").Value = 0.5516531348638054
$ws.Range("G21").Value = 0.3843355958919403
$ws.Range("H21").Value = 0.5278257652593226
$ws.Range("K21").Value = 0.6091354686973887
$ws.Range("L21").Value = 0.1718955583621948
$ws.Range("N21").Value = 1.566367797154538
$ws.Range("O21").Value = 1.766591043415985
$ws.Range("B22").Value = 1.09758679173558
$ws.Range("D22").Value = 0.03906458354920517
$ws.Range("E22").Value = 0.3613880540074916
$ws.Range("F22").Value = 0.5561752365463875
$ws.Range("G22").Value = 0.3871056367822945
$ws.Range("H22").Value = 0.5269005475208388
$ws.Range("K22").Value = 0.6480099219964472
$ws.Range("L22").Value = 0.1815467597092919
$ws.Range("N22").Value = 1.554963057167946
$ws.Range("O22").Value = 1.770686620736598
$ws.Range("B23").Value = 1.083413389619665
$ws.Range("D23").Value = 0.03787004149145901
$ws.Range("E23").Value = 0.3624581150872146
$ws.Range("F23").Value = 0.5537354752168682
$ws.Range("G23").Value = 0.3856062194965091
$ws.Range("H23").Value = 0.527376184402911
$ws.Range("K23").Value = 0.6272766218785648
$ws.Range("L23").Value = 0.1763930416928332
$ws.Range("N23").Value = 1.561006363142965
$ws.Range("O23").Value = 1.768418832410958
$ws.Range("B24").Value = 1.030456053783837
$ws.Range("D24").Value = 0.03332557721726914
$ws.Range("E24").Value = 0.3667938179358501
$ws.Range("F24").Value = 0.5450718275834774
$ws.Range("G24").Value = 0.3803884724759001
$ws.Range("H24").Value = 0.5295770171828593
$ws.Range("K24").Value = 0.5484334590381366
$ws.Range("L24").Value = 0.1569340445048368
$ws.Range("N24").Value = 1.584852636984259
$ws.Range("O24").Value = 1.761629432165591
$ws.Range("B25").Value = 0.9748981431891934
$ws.Range("D25").Value = 0.02838158418925474
$ws.Range("E25").Value = 0.3720761118742963
$ws.Range("F25").Value = 0.5369603321225824
$ws.Range("G25").Value = 0.3757442028136992
$ws.Range("H25").Value = 0.5328060681943327
$ws.Range("K25").Value = 0.4627350178062386
$ws.Range("L25").Value = 0.136077318696934
$ws.Range("N25").Value = 1.602993787077679
$ws.Range("O25").Value = 1.758809453234903

Write-Host "Applied 380 kV case values"